$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.114.04'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.904.98'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').Value = '253.39'
$ws.Range('E5').Value = '  +3.28%  '
$ws.Range('D6').Value = '0.695'
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').Value = '41.39'
$ws.Range('E8').Value = '  +2.50%  '
$ws.Range('D9').Value = '0.359'
$ws.Range('E9').Value = '  +4.36%  '
$ws.Range('D10').Value = '52.51'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('E11').Value = '  +4.90%  '
$ws.Range('D12').Value = '0.0981'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '13.21'
$ws.Range('E13').Value = '  +5.92%  '
$ws.Range('D14').Value = '2.182.20'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('E15').Value = '  +4.86%  '
$ws.Range('D16').Value = '5.01'
$ws.Range('E16').Value = '  +5.40%  '
$ws.Range('D17').Value = '1.914.50'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '35.121.89'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = '73.73'
$ws.Range('E19').Value = '  +2.60%  '
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').Value = '242.99'
$ws.Range('D22').Value = '12.96'
$ws.Range('E22').Value = '  +3.54%  '
$ws.Range('E23').Value = '  +6.26%  '
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '2.43'
$ws.Range('E25').Value = '  +5.53%  '
$ws.Range('D26').Value = '2.31'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').Value = '167.96'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '8.58'
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('D29').Value = '18.55'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('D30').Value = '0.130'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').Value = '4.128.14'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').Value = '0.0604'
$ws.Range('E32').Value = '  +7.67%  '
$ws.Range('E33').Value = '  +10.15%  '
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  +8.85%  '
$ws.Range('D36').Value = '4.22'
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').Value = '0.853'
$ws.Range('E38').Value = '  -4.78%  '
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').Value = '103.68'
$ws.Range('E40').Value = '  +16.27%  '
$ws.Range('D41').Value = '17.31'
$ws.Range('E41').Value = '  +8.40%  '
$ws.Range('D42').Value = '0.0216'
$ws.Range('E42').Value = '  +4.11%  '
$ws.Range('E43').Value = '  +1.83%  '
$ws.Range('D44').Value = '0.0651'
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.307.15'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '12.66'
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '6.59'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').Value = '0.0748'
$ws.Range('E51').Value = '  +6.75%  '
